# Insert a new weekly price-report row for "Plátano" (Terminal Hortofrutícola
# Agro Chillán) right after the existing row 506, shifting every subsequent
# data row (old 507..631) down by one (new 508..632), and populate the newly
# inserted row 507 with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 507:631 down to 508:632, leaving an empty row 507 behind
# (formatting of column D - the date column - is inherited from the row
# above, matching the s="2" date style already used throughout the table).
$ws.Rows("507:507").Insert()

# Fill the newly-opened row 507 with the new data point.
$ws.Range("A507").Value = 7
$ws.Range("B507").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C507").Value = "Ñuble"
$ws.Range("D507").Value = 44754
$ws.Range("E507").Value = 16
$ws.Range("F507").Value = "Fruta"
$ws.Range("G507").Value = 100108
$ws.Range("H507").Value = "Tropicales y subtropicales"
$ws.Range("I507").Value = 100108006
$ws.Range("J507").Value = "Plátano"
$ws.Range("K507").Value = "Sin especificar"
$ws.Range("L507").Value = "Primera Pintón"
$ws.Range("M507").Value = 160
$ws.Range("N507").Value = 23000
$ws.Range("O507").Value = 24000
$ws.Range("P507").Value = 23500
$ws.Range("Q507").Value = "$/caja 20 kilos"
$ws.Range("R507").Value = "Ecuador"
$ws.Range("S507").Value = 1175
$ws.Range("T507").Value = 20
